$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Baz changes" (it is split into
# a "Baz chan" run, a _GoBack bookmark, then a "ges" run).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text -like "Baz chan*") {
        $target = $par
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range

    # Range covering the paragraph's visible text only (exclude the
    # trailing paragraph mark character at End-1).
    $textRange = $d.Range($full.Start, $full.End - 1)

    # Replace the whole run of text ("Baz changes", together with the
    # _GoBack bookmark that sits in the middle of it) with the new
    # wording: "Hi its me abhishek." - "its" is wrapped in the spelling
    # proof-error markers Word inserts for a word it flags, and the
    # _GoBack bookmark is preserved at the end of the new text.
    $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:r><w:t xml:space="preserve">Hi </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>its</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> me abhishek.</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'

    $textRange.InsertXML($newXml)
}
